$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New data blocks mirroring columns C:D, E:F, ... pattern
# Column O/P (header 250 / kg->q 2.5)
$ws.Range("O1").Value = 250
$ws.Range("O2").Value = 2.5
$ws.Range("O4").Value = 1
$ws.Range("P4").Value = 3
$ws.Range("O5").Value = 0.5
$ws.Range("P5").Value = 2.5
$ws.Range("O6").Value = 0.2
$ws.Range("P6").Value = 2.6
$ws.Range("O7").Value = 0.1
$ws.Range("P7").Value = 2.5

# Column Q/R (header 280 / kg->q 2.8)
$ws.Range("Q1").Value = 280
$ws.Range("Q2").Value = 2.8
$ws.Range("Q4").Value = 1
$ws.Range("R4").Value = 3
$ws.Range("Q5").Value = 0.5
$ws.Range("R5").Value = 3
$ws.Range("Q6").Value = 0.2
$ws.Range("R6").Value = 2.8
$ws.Range("Q7").Value = 0.1
$ws.Range("R7").Value = 2.8

# Column S/T (header 290 / kg->q 2.5)
$ws.Range("S1").Value = 290
$ws.Range("S2").Value = 2.5
$ws.Range("S4").Value = 1
$ws.Range("T4").Value = 3
$ws.Range("S5").Value = 0.5
$ws.Range("T5").Value = 3
$ws.Range("S6").Value = 0.2
$ws.Range("T6").Value = 3
$ws.Range("S7").Value = 0.1
$ws.Range("T7").Value = 2.9

# Column U/V (header 191.56 / kg->q 1.9156)
$ws.Range("U1").Value = 191.56
$ws.Range("U2").Value = 1.9156
$ws.Range("U4").Value = 1
$ws.Range("V4").Value = 2
$ws.Range("U5").Value = 0.5
$ws.Range("V5").Value = 2
$ws.Range("U6").Value = 0.2
$ws.Range("V6").Value = 2
$ws.Range("U7").Value = 0.1
$ws.Range("V7").Value = 2

# Copy the border/format styling from the previous block (M:N) onto each new
# pair of columns (copying a 2-col source onto a 1-col target tiles it across
# both destination columns, e.g. O gets M's format and P gets N's format)
$ws.Range("M1:N7").Copy()
$ws.Range("O1:O7").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("Q1:Q7").PasteSpecial(-4122)
$ws.Range("S1:S7").PasteSpecial(-4122)
$ws.Range("U1:U7").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$ws.Range("V7").Select()
